$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "2009年" row). This shifts the "2010年" row (was row 3)
# up to row 2, and the "2011年" row (was row 4) up to row 3.
$ws.Rows.Item(2).Delete()
